$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tiempo extra")
$ws.Range("F2").Formula = "=IF(D2-C2<0,D2-C2,0)"
